$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns per the crypto data refresh.
# D-column values are forced to remain text (NumberFormat "@") to preserve
# exact string formatting (e.g. "307.50", "1.781.87"), then the style is reset
# back to Normal so no formatting/style delta is introduced on the cell.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.732.15"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.82%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.773.06"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -3.17%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.005"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.31%  "

$ws.Range("E5").Value = "  +0.31%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "307.50"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.57%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4386"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.31%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3646"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.33%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07193"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.14%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8375"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.89%  "

$ws.Range("E11").Value = "  -1.84%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.781.87"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -7.14%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.251"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.48%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.345"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.82%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06797"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.01%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.008"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.39%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "79.40"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.52%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008697"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.12%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.004"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.28%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.92"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.93%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "26.622.84"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.96%  "

$ws.Range("E22").Value = "  -2.36%  "

$ws.Range("E23").Value = "  +2.03%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.967.58"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -7.90%  "

$ws.Range("E25").Value = "  -4.35%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "153.24"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.69%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.18"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.63%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.060"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.40%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "114.36"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.18%  "

$ws.Range("E30").Value = "  -9.89%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08998"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.76%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7181"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.78%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.321"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.61%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.799"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -6.27%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.085"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.82%  "

$ws.Range("E36").Value = "  +0.34%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.072"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.35%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05100"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.11%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01888"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.30%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.4923"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.82%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1607"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.05%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.550"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -8.63%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.124"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -6.14%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "7.919"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.41%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "104.80"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.87%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.005"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.41%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.08"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.49%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.06216"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.19%  "

$ws.Range("E49").Value = "  -3.95%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.576"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.79%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.704"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.47%  "
